$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mestre Malte")

# New numeric coefficients added along row 3 (C3:I3, K3), plus a text value in J3
# (shared string "0,10.5" -- not a valid number because of the mixed comma/period,
# so it is stored as a text cell).
$ws.Range("C3").Value = 0.09
$ws.Range("D3").Value = 0.12
$ws.Range("E3").Value = 0.15
$ws.Range("F3").Value = 0.13
$ws.Range("G3").Value = 0.22
$ws.Range("H3").Value = 0.18
$ws.Range("I3").Value = 0.1
$ws.Range("J3").Value = "0,10.5"
$ws.Range("K3").Value = 0.11

# Match the centered formatting used by the rest of row 3 / row 1 for the two
# newly-touched cells that previously had no style at all.
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("K3").HorizontalAlignment = -4108

# Column I widened to fit the new values.
$ws.Columns.Item(9).ColumnWidth = 8.5

# Selection moved from I19 to K8.
$ws.Range("K8").Select()
